# ============================================================================
# "identificacion de tablas.xlsx" update
# Commit: Actualizacion hasta recuperar, registrar, login clase validar datos
#
# Adds an "Obciones" (Opciones) flag column (F) marking which fields are
# required/repeatable, fills in the missing "vachar" characteristic for the
# Nit_Usuario/Clave rows, and appends three new data-dictionary rows:
#   registro (Activo / char), Tipo_Usuario (role list / vachar), equipo (list)
# ============================================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Three new data-dictionary rows: registro, Tipo_Usuario, equipo
#    (written first so new shared-string entries land in the same order as
#    the authored workbook: registro, Tipo_Usuario, Administrador..., Activo)
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = "registro"
$ws.Range("C22").Value = "Tipo_Usuario"
$ws.Range("D22").Value = "Administrador, Jugador, Generador, Representante"
$ws.Range("D21").Value = "Activo"

# ---------------------------------------------------------------------------
# 2. New "Obciones" column header + flags on the existing rows 9-20
# ---------------------------------------------------------------------------
$ws.Range("F9").Value = "Obciones"

$ws.Range("F16").Value = "varias"
$ws.Range("F17").Value = "varias"
$ws.Range("F19").Value = "varias"

$ws.Range("F10").Value = 1
$ws.Range("F11").Value = 1
$ws.Range("F12").Value = 1
$ws.Range("F13").Value = 1
$ws.Range("F14").Value = 1
$ws.Range("F15").Value = 1
$ws.Range("F18").Value = 1
$ws.Range("F20").Value = 1

# A new empty merged pair next to the "Obciones" header column, mirroring the
# existing G9:H9 merged pair.
$ws.Range("H11:I11").Merge()
$ws.Range("G9:H9").Copy()
$ws.Range("H11:I11").PasteSpecial(-4122)
$ws.Range("H11").Value = $null

# ---------------------------------------------------------------------------
# 3. Fill in the missing "vachar" characteristic for Nit_Usuario / Clave
# ---------------------------------------------------------------------------
$ws.Range("E19").Value = "vachar"
$ws.Range("E20").Value = "vachar"

# ---------------------------------------------------------------------------
# 4. Finish the three new rows: flags, "char" characteristic, and the
#    trailing "equipo" row.
# ---------------------------------------------------------------------------
$ws.Range("F21").Value = 1
$ws.Range("E21").Value = "char"

$ws.Range("C23").Value = "equipo"
$ws.Range("D23").Value = "Atleticos(fuerte,regular,novato, aceson-regular, asenso-novato), Space, Unicor, Pastora, "

$ws.Range("E22").Value = "vachar"
$ws.Range("F22").Value = "varias"

# --- Formatting for the new label cells (C21:C23), matching the themed
#     fill used by the rest of column C, but without top/bottom rules since
#     they sit past the original bottom-bordered table edge.
$labels = @("C21", "C22", "C23")
foreach ($addr in $labels) {
    $rng = $ws.Range($addr)
    $rng.Interior.ThemeColor = 10
    $rng.Borders.Item(7).LineStyle = 1
    $rng.Borders.Item(7).Weight = -4138
    $rng.Borders.Item(10).LineStyle = 1
    $rng.Borders.Item(10).Weight = -4138
    $rng.VerticalAlignment = -4108
}

# --- Formatting for the new "caracteristica" cells (E21:E22): lighter tinted
#     fill with only a right-hand rule, matching the rest of column E.
$ws.Range("D11").Copy()
$carac = @("E21", "E22")
foreach ($addr in $carac) {
    $ws.Range($addr).PasteSpecial(-4122)
    $rng = $ws.Range($addr)
    $rng.Borders.Item(7).LineStyle = -4142
    $rng.Borders.Item(8).LineStyle = -4142
    $rng.Borders.Item(9).LineStyle = -4142
    $rng.Borders.Item(10).LineStyle = 1
    $rng.Borders.Item(10).Weight = -4138
}
# restore the values PasteSpecial may have touched (format-only paste keeps
# them, but re-assert to be safe)
$ws.Range("E21").Value = "char"
$ws.Range("E22").Value = "vachar"

# ---------------------------------------------------------------------------
# 4. Move the active selection down to the newly edited area
# ---------------------------------------------------------------------------
$ws.Range("D20").Select()
